$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row appended to the profit log for 2025-11-19.
# Force the date cell to be treated as plain text (matching the existing
# "MM/DD/YYYY" text entries in column A) instead of being auto-parsed into
# a date serial number, then clear the temporary formatting so the cell
# ends up with no explicit style, just like its neighbors.
$ws.Range("A94").NumberFormat = "@"
$ws.Range("A94").Value = "11/19/2025"
$ws.Range("A94").ClearFormats()

$ws.Range("B94").Value = 8184.74
